## Add files via upload
## - adds a new "Posibles ubicaciones" sheet with the catalogue of possible
##   library locations
## - adds a "bibliotecas" defined name pointing at that catalogue
## - adds "Propietario" / "Ubicación" columns (I, J) to the
##   "STATUS ACCOUNTING CIs" sheet, with a list-validation dropdown on J4
##   driven by the new defined name

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("STATUS ACCOUNTING CIs")
$cis = $wb.Worksheets.Item("CIS")

# ---------------------------------------------------------------------
# 1. Add the new "Posibles ubicaciones" sheet right after "CIS"
# ---------------------------------------------------------------------
$locSheet = $wb.Worksheets.Add($null, $cis)
$locSheet.Name = "Posibles ubicaciones"

# ---------------------------------------------------------------------
# 2. Defined name used by the new dropdown validation
# ---------------------------------------------------------------------
$wb.Names.Add("bibliotecas", "='Posibles ubicaciones'!`$C`$6:`$C`$9")

# ---------------------------------------------------------------------
# 3. New header columns on "STATUS ACCOUNTING CIs" (I3 / J3)
#    (written before the catalogue sheet values below so the shared
#    string table grows in the same order as the source workbook)
# ---------------------------------------------------------------------
$ws1.Range("I3").Value = "Propietario (la persona que lo ha creado)"
$ws1.Range("J3").Value = "Ubicación"

# copy the existing header formatting (bold font, centered, thin border)
# onto the two new header cells
$headerFmt = $ws1.Range("C3")
$headerFmt.Copy()
$ws1.Range("I3").PasteSpecial(-4122)
$headerFmt.Copy()
$ws1.Range("J3").PasteSpecial(-4122)
$ws1.Range("I3").WrapText = $true

# copy the plain bordered body-cell formatting onto the new I/J columns
$bodyFmt = $ws1.Range("C4")
$bodyFmt.Copy()
$ws1.Range("I4:J16").PasteSpecial(-4122)

$ws1.Rows.Item(3).RowHeight = 30
$ws1.Columns.Item(9).ColumnWidth = 33.17
$ws1.Columns.Item(10).ColumnWidth = 11.31

# ---------------------------------------------------------------------
# 4. List-validation dropdown on J4 referencing the new defined name
# ---------------------------------------------------------------------
$valCell = $ws1.Range("J4")
$valCell.Validation.Delete()
$valCell.Validation.Add(3, 1, 1, "=bibliotecas")
$valCell.Validation.IgnoreBlank = $true
$valCell.Validation.InCellDropdown = $true
$valCell.Validation.ShowInput = $true
$valCell.Validation.ShowError = $true

[void]$ws1.Range("J4").Select()

# ---------------------------------------------------------------------
# 5. Populate "Posibles ubicaciones" catalogue (C5:C9)
# ---------------------------------------------------------------------
$locSheet.Range("C5").Value = "Posibles Ubicaciones "
$locSheet.Range("C6").Value = "Biblioteca de trabajo"
$locSheet.Range("C7").Value = "Biblioteca de integración"
$locSheet.Range("C8").Value = "Biblioteca de soporte"
$locSheet.Range("C9").Value = "Biblioteca de producción"

$headerFmt.Copy()
$locSheet.Range("C5").PasteSpecial(-4122)
$bodyFmt.Copy()
$locSheet.Range("C6:C9").PasteSpecial(-4122)

$locSheet.Columns.Item(3).ColumnWidth = 23.17
